# Rail freight / rail passenger diesel & electric efficiency trajectories:
# hold the base-year (column J) value flat across all forecast years
# (columns K through AS) for the affected rows, instead of letting the
# value escalate year over year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 5, 6, 9)

foreach ($r in $rows) {
    $baseValue = $ws.Cells.Item($r, 10).Value2   # column J = 10
    $ws.Range($ws.Cells.Item($r, 11), $ws.Cells.Item($r, 45)).Value = $baseValue  # columns K(11) .. AS(45)
}
